$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data for the new rows (3, 4, 5).
# Columns A-F are stored as text, columns G-M as numbers.
$data = @(
    @("120/SUP 2", "Supervision", "O3546845", "KHALID VAVA ", "ds", "mensuelle", 18000, 0, 0, 0, 0, 9000, 18000),
    @("120/SUP 2", "Supervision", "L3659652", "FATIMA FAFA", "ds", "mensuelle", 12000, 0, 10, 600, 0, 6000, 11400),
    @("001/CASA NORD", "Direction régionale", "K6546841", "IMANE FAFA", "ds", "mensuelle", 48000, 36000, 15, 5400, 5400, 12000, 42600)
)

$startRow = 3
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $startRow + $i
    $values = $data[$i]

    # Columns A-F: text values (none of them look purely numeric, so Excel
    # keeps them as plain text without needing a quote-prefix/text style).
    for ($col = 1; $col -le 6; $col++) {
        $cell = $ws.Cells.Item($row, $col)
        $cell.Value = $values[$col - 1]
    }

    # Columns G-M: numeric values
    for ($col = 7; $col -le 13; $col++) {
        $cell = $ws.Cells.Item($row, $col)
        $cell.Value = $values[$col - 1]
    }
}
